# Apply updated dSF (column F) values to Sheet1, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    4  = 2
    5  = -3
    6  = -1
    7  = -1
    8  = -2
    9  = 2
    11 = -2
    12 = 2
    13 = 2
    14 = -1
    15 = 1
    16 = 1
    17 = 1
    18 = -4
    19 = 3
    21 = -4
    22 = -1
    23 = 0
    24 = -1
    25 = 2
    27 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
